$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4389876.5
$ws.Range("I51").Value = 3915.7856
$ws.Range("K51").Value = 3915.7856
$ws.Range("M51").Value = -3431.7856

$ws.Range("H138").Value = 2418.5173
$ws.Range("I138").Value = 2349.2666
$ws.Range("K138").Value = 7047.7998
$ws.Range("M138").Value = -1907.7998

$ws.Range("H141").Value = 5664.75
$ws.Range("J141").Value = 19999
$ws.Range("L141").Value = 59997
$ws.Range("N141").Value = -70357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 9521.200000000001
$ws.Range("I31").Value = 4271.75
$ws.Range("K31").Value = 4271.75
$ws.Range("M31").Value = -3977.75

$ws.Range("H32").Value = 2224.875
$ws.Range("I32").Value = 2208.1409
$ws.Range("K32").Value = 2208.1409
$ws.Range("M32").Value = -1921.1409

$ws.Range("H34").Value = 27766.2
$ws.Range("J34").Value = 30493
$ws.Range("L34").Value = 30493
$ws.Range("N34").Value = -31035

$ws.Range("H45").Value = 5755.1514
$ws.Range("I45").Value = 7258.4
$ws.Range("K45").Value = 7258.4
$ws.Range("M45").Value = -6881.4

$ws.Range("H61").Value = 2908.375
$ws.Range("I61").Value = 2208.8333
$ws.Range("J61").Value = 5007
$ws.Range("K61").Value = 2208.8333
$ws.Range("L61").Value = 5007
$ws.Range("M61").Value = -1996.8333
$ws.Range("N61").Value = -5431

$ws.Range("H74").Value = 3985.238
$ws.Range("J74").Value = 3997.6667
$ws.Range("L74").Value = 3997.6667
$ws.Range("N74").Value = -5745.6667

$ws.Range("H77").Value = 3985.238
$ws.Range("J77").Value = 3997.6667
$ws.Range("L77").Value = 19988.3335
$ws.Range("N77").Value = -28724.3335

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H102").Value = 1902.4667
$ws.Range("I102").Value = 1881.2858
$ws.Range("J102").Value = 2199
$ws.Range("K102").Value = 1881.2858
$ws.Range("L102").Value = 2199
$ws.Range("M102").Value = -259.2858000000001
$ws.Range("N102").Value = -5443

$ws.Range("H110").Value = 2499.5
$ws.Range("J110").Value = 3000
$ws.Range("L110").Value = 3000
$ws.Range("N110").Value = -7090

$ws.Range("H132").Value = 2744.5833
$ws.Range("I132").Value = 2348.889
$ws.Range("K132").Value = 7046.667
$ws.Range("M132").Value = -4516.667

$ws.Range("H136").Value = 2908.375
$ws.Range("I136").Value = 2208.8333
$ws.Range("J136").Value = 5007
$ws.Range("K136").Value = 6626.499899999999
$ws.Range("L136").Value = 15021
$ws.Range("M136").Value = -4076.499899999999
$ws.Range("N136").Value = -20121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 79326.46000000001
$ws.Range("I107").Value = 85602.586
$ws.Range("J107").Value = 4013
$ws.Range("K107").Value = 85602.586
$ws.Range("L107").Value = 4013
$ws.Range("M107").Value = -83682.586
$ws.Range("N107").Value = -7853

$ws.Range("H134").Value = 121387.52
$ws.Range("I134").Value = 177469.17
$ws.Range("J134").Value = 2214
$ws.Range("K134").Value = 532407.51
$ws.Range("L134").Value = 6642
$ws.Range("M134").Value = -529872.51
$ws.Range("N134").Value = -11712

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1922.4
$ws.Range("I94").Value = 1880.4
$ws.Range("J94").Value = 1964.4
$ws.Range("K94").Value = 1880.4
$ws.Range("L94").Value = 1964.4
$ws.Range("M94").Value = -1429.4
$ws.Range("N94").Value = -2866.4

$ws.Range("H99").Value = 19811.875
$ws.Range("I99").Value = 28899.6
$ws.Range("K99").Value = 28899.6
$ws.Range("M99").Value = -27401.6

$ws.Range("H126").Value = 19811.875
$ws.Range("I126").Value = 28899.6
$ws.Range("K126").Value = 86698.79999999999
$ws.Range("M126").Value = -84228.79999999999

$ws.Range("H141").Value = 388850.44
$ws.Range("J141").Value = 388850.44
$ws.Range("L141").Value = 388850.44
$ws.Range("N141").Value = -399210.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 499.6
$ws.Range("I2").Value = 716.4
$ws.Range("J2").Value = 66
$ws.Range("K2").Value = 4298.4
$ws.Range("L2").Value = 396
$ws.Range("M2").Value = -4185.4
$ws.Range("N2").Value = -622

$ws.Range("H97").Value = 949
$ws.Range("I97").Value = 997
$ws.Range("J97").Value = 925
$ws.Range("K97").Value = 2991
$ws.Range("L97").Value = 2775
$ws.Range("M97").Value = -2495
$ws.Range("N97").Value = -3767

$ws.Range("H107").Value = 599.75
$ws.Range("I107").Value = 1000
$ws.Range("J107").Value = 199.5
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 598.5
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -4438.5

$ws.Range("H116").Value = 5499
$ws.Range("I116").Value = 2498.3333
$ws.Range("K116").Value = 7494.999899999999
$ws.Range("M116").Value = -4052.999899999999

$ws.Range("H122").Value = 766.36365
$ws.Range("I122").Value = 698.2857
$ws.Range("K122").Value = 6284.571300000001
$ws.Range("M122").Value = -3834.571300000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3087.8262
$ws.Range("I102").Value = 2989.0527
$ws.Range("J102").Value = 3557
$ws.Range("K102").Value = 2989.0527
$ws.Range("L102").Value = 3557
$ws.Range("M102").Value = -1367.0527
$ws.Range("N102").Value = -6801

$ws.Range("H132").Value = 7013.931
$ws.Range("I132").Value = 6440.1333
$ws.Range("J132").Value = 7628.7144
$ws.Range("K132").Value = 19320.3999
$ws.Range("L132").Value = 22886.1432
$ws.Range("M132").Value = -16790.3999
$ws.Range("N132").Value = -27946.1432

$ws.Range("H136").Value = 33962.668
$ws.Range("J136").Value = 33962.668
$ws.Range("L136").Value = 101888.004
$ws.Range("N136").Value = -106988.004

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 23085
$ws.Range("J63").Value = 23085
$ws.Range("L63").Value = 23085
$ws.Range("N63").Value = -24583

$ws.Range("H66").Value = 23085
$ws.Range("J66").Value = 23085
$ws.Range("L66").Value = 69255
$ws.Range("N66").Value = -76743

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3967.818
$ws.Range("I96").Value = 4774.8335
$ws.Range("K96").Value = 4774.8335
$ws.Range("M96").Value = -3401.8335
